# Applies the crypto price/volume refresh described in the commit diff.
# Cell values are plain text (inline strings) in the source workbook, so any
# price that happens to look like a number (e.g. "242.82") is written with the
# cell pre-formatted as Text ("@") to stop Excel from silently converting it to
# a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.160.74'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '1.829.08'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.82'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6199'
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07344'
$ws.Range("E8").Value = '  -1.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2904'
$ws.Range("E9").Value = '  -0.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.25'
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07630'
$ws.Range("E11").Value = '  -0.85%  '
$ws.Range("D12").Value = '1.840.39'
$ws.Range("E12").Value = '  +0.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.960'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6689'
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("E15").Value = '  -0.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008977'
$ws.Range("E16").Value = '  -1.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.833'
$ws.Range("E17").Value = '  -1.34%  '
$ws.Range("D18").Value = '29.152.93'
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("D19").Value = '2.085.46'
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '235.83'
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("E21").Value = '  -1.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.343'
$ws.Range("E23").Value = '  +1.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.50'
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("E26").Value = '  -0.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.518'
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("E28").Value = '  -1.24%  '
$ws.Range("E29").Value = '  -0.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05834'
$ws.Range("E30").Value = '  +6.09%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.221'
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.075'
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.084'
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.855'
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7245'
$ws.Range("E36").Value = '  -1.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.607'
$ws.Range("E37").Value = '  -2.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.849'
$ws.Range("E38").Value = '  +2.57%  '
$ws.Range("D39").Value = '1.227.30'
$ws.Range("E39").Value = '  +1.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01759'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.222'
$ws.Range("E41").Value = '  -3.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9077'
$ws.Range("E42").Value = '  +1.98%  '
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.80'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '1.987.97'
$ws.Range("E45").Value = '  +0.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.76'
$ws.Range("E46").Value = '  +0.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5038'
$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4033'
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.117'
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("E50").Value = '  -2.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1128'
$ws.Range("E51").Value = '  +2.95%  '
